# Auto-generated edit script to update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '59.199.72'
$ws.Cells.Item(2, 5).Value = '  -6.07%  '

$ws.Cells.Item(3, 4).Value = '2.453.97'
$ws.Cells.Item(3, 5).Value = '  -8.53%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '539.01'
$ws.Cells.Item(5, 5).Value = '  -2.98%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '147.53'
$ws.Cells.Item(6, 5).Value = '  -7.10%  '

$ws.Cells.Item(7, 5).Value = '  -0.19%  '

$ws.Cells.Item(8, 5).Value = '  -4.24%  '

$ws.Cells.Item(9, 4).Value = '2.472.33'
$ws.Cells.Item(9, 5).Value = '  -7.99%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0993'
$ws.Cells.Item(10, 5).Value = '  -6.35%  '

$ws.Cells.Item(11, 5).Value = '  -2.54%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '5.30'

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.352'
$ws.Cells.Item(13, 5).Value = '  -4.51%  '

$ws.Cells.Item(14, 4).Value = '2.895.00'
$ws.Cells.Item(14, 5).Value = '  -8.31%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '24.08'
$ws.Cells.Item(15, 5).Value = '  -8.76%  '

$ws.Cells.Item(16, 4).Value = '59.108.62'
$ws.Cells.Item(16, 5).Value = '  -6.02%  '

$ws.Cells.Item(17, 5).Value = '  -6.28%  '

$ws.Cells.Item(18, 4).Value = '2.522.22'
$ws.Cells.Item(18, 5).Value = '  -6.02%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '11.13'
$ws.Cells.Item(19, 5).Value = '  -6.80%  '

$ws.Cells.Item(20, 5).Value = '  -6.26%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '324.62'

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.967'
$ws.Cells.Item(22, 5).Value = '  -3.28%  '

$ws.Cells.Item(23, 5).Value = '  -9.06%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.460'
$ws.Cells.Item(24, 5).Value = '  -10.02%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '60.67'
$ws.Cells.Item(25, 5).Value = '  -4.27%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.161'
$ws.Cells.Item(26, 5).Value = '  -4.31%  '

$ws.Cells.Item(27, 5).Value = '  -2.04%  '

$ws.Cells.Item(29, 5).Value = '  -6.56%  '

$ws.Cells.Item(30, 5).Value = '  -7.98%  '

$ws.Cells.Item(31, 2).Value = 'PEPE'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(31, 4).Value = '0.0₃0775'
$ws.Cells.Item(31, 5).Value = '  -9.92%  '

$ws.Cells.Item(32, 2).Value = 'Fetch.AI'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.26'
$ws.Cells.Item(32, 5).Value = '  -8.09%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.997'
$ws.Cells.Item(33, 5).Value = '  -0.11%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '158.43'
$ws.Cells.Item(34, 5).Value = '  -4.29%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '4.55'
$ws.Cells.Item(35, 5).Value = '  -7.20%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.38'
$ws.Cells.Item(36, 5).Value = '  -7.61%  '

$ws.Cells.Item(37, 5).Value = '  -5.91%  '

$ws.Cells.Item(38, 5).Value = '  -2.08%  '

$ws.Cells.Item(39, 2).Value = 'Bittensor'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '319.87'
$ws.Cells.Item(39, 5).Value = '  -8.48%  '

$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '5.90'
$ws.Cells.Item(40, 5).Value = '  -7.67%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '36.63'
$ws.Cells.Item(41, 5).Value = '  -4.33%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.838'
$ws.Cells.Item(42, 5).Value = '  -12.48%  '

$ws.Cells.Item(43, 5).Value = '  -7.89%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.996'
$ws.Cells.Item(44, 5).Value = '  -0.26%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '10.72'
$ws.Cells.Item(45, 5).Value = '  -2.78%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.585'
$ws.Cells.Item(46, 5).Value = '  -5.45%  '

$ws.Cells.Item(47, 5).Value = '  -3.43%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0526'
$ws.Cells.Item(48, 5).Value = '  -6.52%  '

$ws.Cells.Item(49, 2).Value = 'VeChain'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0229'
$ws.Cells.Item(49, 5).Value = '  -5.10%  '

$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '18.57'
$ws.Cells.Item(50, 5).Value = '  -9.11%  '

$ws.Cells.Item(51, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '19.00'
$ws.Cells.Item(51, 5).Value = '  -8.96%  '
